$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.004257666666666667
$ws.Range("H2").Value = 0.012773
$ws.Range("M2").Value = 31.57670733333334
$ws.Range("N2").Value = 94.73012200000001
$ws.Range("O2").Value = 0.3460135388827164
$ws.Range("P2").Value = 0.3460135388827164
$ws.Range("Q2").Value = 0.1344430942562222
$ws.Range("R2").Value = 1.209987848306
$ws.Range("S2").Value = 0.3460135388827164
$ws.Range("T2").Value = 0.3460135388827164

# Row 3
$ws.Range("G3").Value = 0.004257666666666667
$ws.Range("H3").Value = 0.012773
$ws.Range("O3").Value = 0.4506220157900242
$ws.Range("P3").Value = 0.4506220157900241
$ws.Range("Q3").Value = 0.1750885769915555
$ws.Range("R3").Value = 1.575797192924
$ws.Range("S3").Value = 0.4506220157900242
$ws.Range("T3").Value = 0.4506220157900241

# Row 4
$ws.Range("G4").Value = 0.004257666666666667
$ws.Range("H4").Value = 0.012773
$ws.Range("M4").Value = 18.55875233333333
$ws.Range("N4").Value = 55.676257
$ws.Range("O4").Value = 0.2033644453272594
$ws.Range("P4").Value = 0.2033644453272593
$ws.Range("Q4").Value = 0.07901698118455555
$ws.Range("R4").Value = 0.711152830661
$ws.Range("S4").Value = 0.2033644453272594
$ws.Range("T4").Value = 0.2033644453272593
